# Update the "Corr / Total" marks in the marksheet summary block.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 "Marking" -> Right column (B11): 3 -> 5
$ws.Range("B11").Value = 5

# Row 12 "Total" -> Right column (B12): 54 -> 90
$ws.Range("B12").Value = 90

# Row 12 "Total" -> Max column (E12): "54/84" -> "90/140"
$ws.Range("E12").Value = "90/140"
